# Closing and opening balance change
#
# G column = closing_balance, H column = opening_balance (see row 1 headers).
# For every data row 2..19 these were per-row numeric literals; the edit
# replaces them with the same two text values ("2000 : 1" / "500 : 1") for
# every row, formatted with wrap-text (matching the style already used by
# G2/H2 and the rest of the formatted columns on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$balanceRange = $ws.Range("G2:H19")
$balanceRange.WrapText = $true

$ws.Range("G2:G19").Value = "2000 : 1"
$ws.Range("H2:H19").Value = "500 : 1"

# Move the selection/active cell to match the new focus area.
[void]$ws.Range("G2:H19").Select()
